$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Germany Landesliga")

$ws.Cells.Item(2, 2).Value = 6760228
$ws.Cells.Item(2, 8).Value = 4
$ws.Cells.Item(2, 9).Value = 1
$ws.Cells.Item(2, 11).Value = 2.25
$ws.Cells.Item(2, 13).Value = 2.625
$ws.Cells.Item(2, 14).Value = 1.75
$ws.Cells.Item(2, 15).Value = 3.6
$ws.Cells.Item(2, 16).Value = 3.8
$ws.Cells.Item(2, 17).Value = -0.5
$ws.Cells.Item(2, 18).Value = 1.8
$ws.Cells.Item(2, 19).Value = 2
$ws.Cells.Item(2, 20).Value = 2.75
$ws.Cells.Item(2, 21).Value = 1.8
$ws.Cells.Item(2, 22).Value = 2
$ws.Cells.Item(2, 23).Value = 0.75
$ws.Cells.Item(2, 26).Value = 0.8
$ws.Cells.Item(2, 28).Value = 0.8
$ws.Cells.Item(2, 29).Value = -1
$ws.Cells.Item(3, 2).Value = 6757276
$ws.Cells.Item(3, 8).Value = 3
$ws.Cells.Item(3, 9).Value = 0
$ws.Cells.Item(3, 11).Value = 2.6
$ws.Cells.Item(3, 13).Value = 2.25
$ws.Cells.Item(3, 14).Value = 2.6
$ws.Cells.Item(3, 15).Value = 3.5
$ws.Cells.Item(3, 16).Value = 2.25
$ws.Cells.Item(3, 17).Value = 0.25
$ws.Cells.Item(3, 18).Value = 1.75
$ws.Cells.Item(3, 19).Value = 2.05
$ws.Cells.Item(3, 20).Value = 3.25
$ws.Cells.Item(3, 21).Value = 1.775
$ws.Cells.Item(3, 22).Value = 2.025
$ws.Cells.Item(3, 23).Value = 1.6
$ws.Cells.Item(3, 26).Value = 0.75
$ws.Cells.Item(3, 28).Value = -0.5
$ws.Cells.Item(3, 29).Value = 0.5125
$ws.Cells.Item(4, 2).Value = 6781315
$ws.Cells.Item(4, 7).Value = "SpVg Porz 1919"
$ws.Cells.Item(4, 8).Value = 2
$ws.Cells.Item(4, 9).Value = 0
$ws.Cells.Item(4, 10).Value = "H"
$ws.Cells.Item(4, 11).Value = 1.25
$ws.Cells.Item(4, 12).Value = 4.75
$ws.Cells.Item(4, 13).Value = 10
$ws.Cells.Item(4, 14).Value = 1.222
$ws.Cells.Item(4, 15).Value = 5.25
$ws.Cells.Item(4, 16).Value = 8.5
$ws.Cells.Item(4, 17).Value = -2
$ws.Cells.Item(4, 18).Value = 1.925
$ws.Cells.Item(4, 19).Value = 1.875
$ws.Cells.Item(4, 20).Value = 3.5
$ws.Cells.Item(4, 21).Value = 1.775
$ws.Cells.Item(4, 22).Value = 1.925
$ws.Cells.Item(4, 23).Value = 0.222
$ws.Cells.Item(4, 25).Value = -1
$ws.Cells.Item(4, 26).Value = 0
$ws.Cells.Item(4, 27).Value = -0
$ws.Cells.Item(4, 29).Value = 0.925
$ws.Cells.Item(5, 2).Value = 6781316
$ws.Cells.Item(5, 7).Value = "SV 09 Arnstadt"
$ws.Cells.Item(5, 8).Value = 0
$ws.Cells.Item(5, 9).Value = 2
$ws.Cells.Item(5, 10).Value = "A"
$ws.Cells.Item(5, 11).Value = 5
$ws.Cells.Item(5, 12).Value = 4.5
$ws.Cells.Item(5, 13).Value = 1.45
$ws.Cells.Item(5, 14).Value = 6.5
$ws.Cells.Item(5, 15).Value = 4.333
$ws.Cells.Item(5, 16).Value = 1.363
$ws.Cells.Item(5, 17).Value = 1.5
$ws.Cells.Item(5, 18).Value = 1.825
$ws.Cells.Item(5, 19).Value = 1.975
$ws.Cells.Item(5, 20).Value = 3
$ws.Cells.Item(5, 21).Value = 1.825
$ws.Cells.Item(5, 22).Value = 1.975
$ws.Cells.Item(5, 23).Value = -1
$ws.Cells.Item(5, 25).Value = 0.363
$ws.Cells.Item(5, 26).Value = -1
$ws.Cells.Item(5, 27).Value = 0.9750000000000001
$ws.Cells.Item(5, 29).Value = 0.9750000000000001
$ws.Cells.Item(12, 2).Value = 7035047
$ws.Cells.Item(12, 8).Value = 3
$ws.Cells.Item(12, 9).Value = 4
$ws.Cells.Item(12, 11).Value = 1.909
$ws.Cells.Item(12, 12).Value = 3.75
$ws.Cells.Item(12, 13).Value = 3.1
$ws.Cells.Item(12, 14).Value = 2.2
$ws.Cells.Item(12, 16).Value = 2.625
$ws.Cells.Item(12, 18).Value = 2
$ws.Cells.Item(12, 19).Value = 1.8
$ws.Cells.Item(12, 20).Value = 3
$ws.Cells.Item(12, 21).Value = 1.825
$ws.Cells.Item(12, 22).Value = 1.975
$ws.Cells.Item(12, 25).Value = 1.625
$ws.Cells.Item(12, 27).Value = 0.8
$ws.Cells.Item(12, 28).Value = 0.825
$ws.Cells.Item(12, 29).Value = -1
$ws.Cells.Item(13, 2).Value = 7035046
$ws.Cells.Item(13, 8).Value = 0
$ws.Cells.Item(13, 9).Value = 2
$ws.Cells.Item(13, 11).Value = 2
$ws.Cells.Item(13, 12).Value = 3.6
$ws.Cells.Item(13, 13).Value = 3
$ws.Cells.Item(13, 14).Value = 2
$ws.Cells.Item(13, 16).Value = 3
$ws.Cells.Item(13, 18).Value = 1.8
$ws.Cells.Item(13, 19).Value = 2
$ws.Cells.Item(13, 20).Value = 2.75
$ws.Cells.Item(13, 21).Value = 1.8
$ws.Cells.Item(13, 22).Value = 2
$ws.Cells.Item(13, 25).Value = 2
$ws.Cells.Item(13, 27).Value = 1
$ws.Cells.Item(13, 28).Value = -1
$ws.Cells.Item(13, 29).Value = 1
$ws.Cells.Item(21, 6).Value = "Cronenberger SC"
$ws.Cells.Item(25, 7).Value = "Cronenberger SC"
$ws.Cells.Item(39, 6).Value = "SSV Markranstadt"
$ws.Cells.Item(50, 7).Value = "FSV Duisburg"
$ws.Cells.Item(53, 7).Value = "FSV Duisburg"
$ws.Cells.Item(55, 6).Value = "SC Dsseldorf West"
$ws.Cells.Item(55, 7).Value = "FC Viersen"
$ws.Cells.Item(61, 6).Value = "SpVg Porz 1919"
$ws.Cells.Item(65, 6).Value = "SV Schott Jena"
$ws.Cells.Item(69, 6).Value = "Cronenberger SC"
$ws.Cells.Item(70, 7).Value = "VfL Viktoria JuchenGarzweiler"
$ws.Cells.Item(71, 7).Value = "Cronenberger SC"
$ws.Cells.Item(83, 2).Value = 8075296
$ws.Cells.Item(83, 6).Value = "FC Monheim"
$ws.Cells.Item(83, 7).Value = "VFB Hilden II"
$ws.Cells.Item(83, 8).Value = 1
$ws.Cells.Item(83, 9).Value = 2
$ws.Cells.Item(83, 10).Value = "A"
$ws.Cells.Item(83, 11).Value = 1.533
$ws.Cells.Item(83, 12).Value = 4.75
$ws.Cells.Item(83, 13).Value = 4
$ws.Cells.Item(83, 14).Value = 1.4
$ws.Cells.Item(83, 15).Value = 5.25
$ws.Cells.Item(83, 16).Value = 5
$ws.Cells.Item(83, 17).Value = -1.5
$ws.Cells.Item(83, 18).Value = 1.975
$ws.Cells.Item(83, 19).Value = 1.825
$ws.Cells.Item(83, 20).Value = 3.75
$ws.Cells.Item(83, 21).Value = 1.9
$ws.Cells.Item(83, 22).Value = 1.9
$ws.Cells.Item(83, 23).Value = -1
$ws.Cells.Item(83, 25).Value = 4
$ws.Cells.Item(83, 26).Value = -1
$ws.Cells.Item(83, 27).Value = 0.825
$ws.Cells.Item(83, 29).Value = 0.8999999999999999
$ws.Cells.Item(84, 2).Value = 8075530
$ws.Cells.Item(84, 6).Value = "TuRU Dsseldorf"
$ws.Cells.Item(84, 7).Value = "DV Solingen"
$ws.Cells.Item(84, 8).Value = 2
$ws.Cells.Item(84, 9).Value = 0
$ws.Cells.Item(84, 10).Value = "H"
$ws.Cells.Item(84, 11).Value = 2.1
$ws.Cells.Item(84, 12).Value = 3.75
$ws.Cells.Item(84, 13).Value = 2.7
$ws.Cells.Item(84, 14).Value = 2.375
$ws.Cells.Item(84, 15).Value = 3.75
$ws.Cells.Item(84, 16).Value = 2.45
$ws.Cells.Item(84, 17).Value = 0
$ws.Cells.Item(84, 18).Value = 1.85
$ws.Cells.Item(84, 19).Value = 1.95
$ws.Cells.Item(84, 20).Value = 3
$ws.Cells.Item(84, 21).Value = 1.85
$ws.Cells.Item(84, 22).Value = 1.95
$ws.Cells.Item(84, 23).Value = 1.375
$ws.Cells.Item(84, 25).Value = -1
$ws.Cells.Item(84, 26).Value = 0.8500000000000001
$ws.Cells.Item(84, 27).Value = -1
$ws.Cells.Item(84, 29).Value = 0.95
$ws.Cells.Item(85, 2).Value = 8075670
$ws.Cells.Item(85, 7).Value = "VfB Frohnhausen"
$ws.Cells.Item(85, 8).Value = 2
$ws.Cells.Item(85, 9).Value = 1
$ws.Cells.Item(85, 10).Value = "H"
$ws.Cells.Item(85, 11).Value = 1.222
$ws.Cells.Item(85, 12).Value = 6.5
$ws.Cells.Item(85, 13).Value = 9.5
$ws.Cells.Item(85, 14).Value = 1.181
$ws.Cells.Item(85, 15).Value = 7
$ws.Cells.Item(85, 16).Value = 11
$ws.Cells.Item(85, 17).Value = -2.5
$ws.Cells.Item(85, 18).Value = 1.925
$ws.Cells.Item(85, 19).Value = 1.875
$ws.Cells.Item(85, 20).Value = 4.25
$ws.Cells.Item(85, 21).Value = 1.8
$ws.Cells.Item(85, 22).Value = 2
$ws.Cells.Item(85, 23).Value = 0.181
$ws.Cells.Item(85, 25).Value = -1
$ws.Cells.Item(85, 27).Value = 0.875
$ws.Cells.Item(85, 29).Value = 1
$ws.Cells.Item(86, 2).Value = 8075593
$ws.Cells.Item(86, 7).Value = "ESC Rellinghausen"
$ws.Cells.Item(86, 8).Value = 0
$ws.Cells.Item(86, 9).Value = 2
$ws.Cells.Item(86, 10).Value = "A"
$ws.Cells.Item(86, 11).Value = 6.5
$ws.Cells.Item(86, 12).Value = 4.5
$ws.Cells.Item(86, 13).Value = 1.363
$ws.Cells.Item(86, 14).Value = 4.75
$ws.Cells.Item(86, 15).Value = 4.5
$ws.Cells.Item(86, 16).Value = 1.5
$ws.Cells.Item(86, 17).Value = 1.25
$ws.Cells.Item(86, 18).Value = 1.875
$ws.Cells.Item(86, 19).Value = 1.925
$ws.Cells.Item(86, 20).Value = 3.75
$ws.Cells.Item(86, 21).Value = 2
$ws.Cells.Item(86, 22).Value = 1.8
$ws.Cells.Item(86, 23).Value = -1
$ws.Cells.Item(86, 25).Value = 0.5
$ws.Cells.Item(86, 27).Value = 0.925
$ws.Cells.Item(86, 29).Value = 0.8
$ws.Cells.Item(89, 2).Value = 8076477
$ws.Cells.Item(89, 8).Value = 2
$ws.Cells.Item(89, 9).Value = 2
$ws.Cells.Item(89, 10).Value = "D"
$ws.Cells.Item(89, 11).Value = 1.833
$ws.Cells.Item(89, 12).Value = 4
$ws.Cells.Item(89, 13).Value = 3.1
$ws.Cells.Item(89, 14).Value = 1.833
$ws.Cells.Item(89, 15).Value = 4
$ws.Cells.Item(89, 16).Value = 3.1
$ws.Cells.Item(89, 17).Value = -0.5
$ws.Cells.Item(89, 18).Value = 1.875
$ws.Cells.Item(89, 19).Value = 1.925
$ws.Cells.Item(89, 20).Value = 3.75
$ws.Cells.Item(89, 21).Value = 1.975
$ws.Cells.Item(89, 22).Value = 1.825
$ws.Cells.Item(89, 23).Value = -1
$ws.Cells.Item(89, 24).Value = 3
$ws.Cells.Item(89, 27).Value = 0.925
$ws.Cells.Item(89, 28).Value = 0.4875
$ws.Cells.Item(89, 29).Value = -0.5
$ws.Cells.Item(90, 2).Value = 8077795
$ws.Cells.Item(90, 8).Value = 5
$ws.Cells.Item(90, 9).Value = 3
$ws.Cells.Item(90, 10).Value = "H"
$ws.Cells.Item(90, 11).Value = 1.142
$ws.Cells.Item(90, 12).Value = 7
$ws.Cells.Item(90, 13).Value = 10
$ws.Cells.Item(90, 14).Value = 1.083
$ws.Cells.Item(90, 15).Value = 11
$ws.Cells.Item(90, 16).Value = 19
$ws.Cells.Item(90, 17).Value = -3.5
$ws.Cells.Item(90, 18).Value = 1.975
$ws.Cells.Item(90, 19).Value = 1.825
$ws.Cells.Item(90, 20).Value = 5
$ws.Cells.Item(90, 21).Value = 1.825
$ws.Cells.Item(90, 22).Value = 1.975
$ws.Cells.Item(90, 23).Value = 0.08299999999999996
$ws.Cells.Item(90, 24).Value = -1
$ws.Cells.Item(90, 27).Value = 0.825
$ws.Cells.Item(90, 28).Value = 0.825
$ws.Cells.Item(90, 29).Value = -1
$ws.Cells.Item(91, 6).Value = "SC Dsseldorf West"
